$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("FUSELAGE")
$ws.Range("A8").Value = "JENKINSON"
$ws.Range("C8").Value = 21031.0
$ws.Range("D8").Value = 211.1876513474806
$ws.Range("A9").Value = "SADRAEY"
$ws.Range("C9").Value = 6396.0
$ws.Range("D9").Value = -5.36083790506938
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 6850.0
$ws.Range("D10").Value = 1.3568261961029933
$ws.Range("A11").Value = "NICOLAI_1984"
$ws.Range("C11").Value = 10430.0
$ws.Range("D11").Value = 54.32871492340938
$ws.Range("A13").Value = "ROSKAM"
$ws.Range("C13").Value = 15196.0
$ws.Range("D13").Value = 124.84939136875636
$ws.Range("A14").Value = "RAYMER"
$ws.Range("C14").Value = 6652.0
$ws.Range("D14").Value = -1.5729039625580858
$ws.Range("A15").Value = "TORENBEEK_1976"
$ws.Range("C15").Value = 10802.0
$ws.Range("D15").Value = 59.8330564336211

$ws = $wb.Worksheets("WING")
$ws.Range("A8").Value = "KROO"
$ws.Range("C8").Value = 7503.0
$ws.Range("D8").Value = 4.734921766886217
$ws.Range("A9").Value = "TORENBEEK_2013"
$ws.Range("C9").Value = 6097.0
$ws.Range("D9").Value = -14.891534317912134
$ws.Range("A10").Value = "TORENBEEK_1982"
$ws.Range("C10").Value = 6272.0
$ws.Range("D10").Value = -12.448696611767248
$ws.Range("A11").Value = "RAYMER"
$ws.Range("C11").Value = 8551.0
$ws.Range("D11").Value = 19.364029858542455

$ws = $wb.Worksheets("HORIZONTAL TAIL")
$ws.Range("A8").Value = "HOWE"
$ws.Range("C8").Value = 1415.0
$ws.Range("D8").Value = 82.06272176132752
$ws.Range("A9").Value = "JENKINSON"
$ws.Range("C9").Value = 700.0
$ws.Range("D9").Value = -9.933635877788502
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 739.0
$ws.Range("D10").Value = -4.915652733836718
$ws.Range("A11").Value = "SADRAEY"
$ws.Range("C11").Value = 1040.0
$ws.Range("D11").Value = 33.81288383871423
$ws.Range("A12").Value = "NICOLAI_2013"
$ws.Range("C12").Value = 416.0
$ws.Range("D12").Value = -46.47484646451431
$ws.Range("A13").Value = "TORENBEEK_1976"
$ws.Range("C13").Value = 52.0
$ws.Range("D13").Value = -93.30935580806428
$ws.Range("A14").Value = "RAYMER"
$ws.Range("C14").Value = 526.0
$ws.Range("D14").Value = -32.32156067388107
$ws.Range("A15").Value = "ROSKAM"
$ws.Range("C15").Value = 1523.0
$ws.Range("D15").Value = 95.95867508304016

$ws = $wb.Worksheets("VERTICAL TAIL")
$ws.Range("A8").Value = "HOWE"
$ws.Range("C8").Value = 1145.0
$ws.Range("D8").Value = 47.322838457045954
$ws.Range("A9").Value = "JENKINSON"
$ws.Range("C9").Value = 502.0
$ws.Range("D9").Value = -35.40955030092832
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 498.0
$ws.Range("D10").Value = -35.924215238769534
$ws.Range("A11").Value = "SADRAEY"
$ws.Range("C11").Value = 749.0
$ws.Range("D11").Value = -3.6289903892336968
$ws.Range("A12").Value = "TORENBEEK_1976"
$ws.Range("C12").Value = 124.0
$ws.Range("D12").Value = -84.04538692692253
$ws.Range("A13").Value = "RAYMER"
$ws.Range("C13").Value = 186.0
$ws.Range("D13").Value = -76.0680803903838
$ws.Range("A14").Value = "ROSKAM"
$ws.Range("C14").Value = 1523.0
$ws.Range("D14").Value = 95.95867508304016

$ws = $wb.Worksheets("NACELLES")
$ws.Range("A10").Value = "JENKINSON"
$ws.Range("C10").Value = 705.0
$ws.Range("D10").Value = 9.806473251252601
$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 687.0
$ws.Range("D12").Value = 7.002903721433387
$ws.Range("A17").Value = "JENKINSON"
$ws.Range("C17").Value = 705.0
$ws.Range("D17").Value = 9.806473251252601
$ws.Range("A19").Value = "ROSKAM"
$ws.Range("C19").Value = 687.0
$ws.Range("D19").Value = 7.002903721433387

$ws = $wb.Worksheets("POWER PLANT")
$ws.Range("A13").Value = "TORENBEEK_1976"
$ws.Range("C13").Value = 2954.0
$ws.Range("D13").Value = 5.323389177130686
$ws.Range("A20").Value = "TORENBEEK_1976"
$ws.Range("C20").Value = 2954.0
$ws.Range("D20").Value = 5.323389177130686

$ws = $wb.Worksheets("LANDING GEARS")
$ws.Range("A9").Value = "TORENBEEK_1976"
$ws.Range("C9").Value = 2683.5099483972144
$ws.Range("D9").Value = -3.1539784813210665
$ws.Range("A11").Value = "TORENBEEK_1976"
$ws.Range("C11").Value = 410.49191368324796
$ws.Range("A13").Value = "TORENBEEK_1976"
$ws.Range("C13").Value = 2273.0180347139667

